# The source change is a SharePoint "document management" custom-XML
# re-shuffle: the deck carries two metadata parts under customXml/ -
#   1) the document-management property bag
#        p:properties > documentManagement > Status ("Not started"),
#        MediaServiceKeyPoints
#   2) the content-type schema
#        ct:contentTypeSchema (...)
# and these two parts swap physical slots: whichever one used to be
# item2.xml/itemProps2.xml becomes item3.xml/itemProps3.xml and vice
# versa. The XML payloads themselves are byte-for-byte unchanged - only
# which numbered part holds which payload differs (this is what Office
# does internally when it re-persists a document carrying SharePoint
# content-organizer metadata).
#
# Reproduce it through CustomXMLParts: pull each part's XML back out,
# delete both parts, then re-Add them in the opposite order so the
# package writer hands the part names out swapped.

function Swap-DocumentManagementCustomXmlParts {
    param($Presentation)

    $parts = $Presentation.CustomXMLParts
    if ($parts -eq $null) { return }

    $count = 0
    try { $count = $parts.Count } catch { $count = 0 }
    if ($count -le 0) { return }

    $propsPart = $null
    $schemaPart = $null

    for ($i = 1; $i -le $count; $i++) {
        $part = $null
        try { $part = $parts.Item($i) } catch { $part = $null }
        if ($part -eq $null) { continue }

        $xml = $null
        try { $xml = $part.XML } catch { $xml = $null }
        if ([string]::IsNullOrEmpty($xml)) { continue }

        if ($xml -like "*p:properties*documentManagement*") {
            $propsPart = $part
        } elseif ($xml -like "*ct:contentTypeSchema*") {
            $schemaPart = $part
        }
    }

    if ($propsPart -ne $null -and $schemaPart -ne $null) {
        $propsXml = $propsPart.XML
        $schemaXml = $schemaPart.XML

        $propsPart.Delete()
        $schemaPart.Delete()

        # Re-insert swapped: the properties payload now lands in the slot
        # that used to hold the schema (item2.xml) and the schema payload
        # lands where the properties used to be (item3.xml).
        [void]$parts.Add($propsXml)
        [void]$parts.Add($schemaXml)
    }
}

$p = $ppt.ActivePresentation

try {
    Swap-DocumentManagementCustomXmlParts -Presentation $p
} catch {
    # Non-fatal: leave the deck untouched if custom XML parts aren't
    # reachable in this environment.
}
